$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert new "MenuStrip" sheet right after "messageBox"
# ------------------------------------------------------------------
$msgBox = $wb.Worksheets.Item("messageBox")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $msgBox)
$ws.Name = "MenuStrip"

# Header row (bold black text, like the other control sheets)
$ws.Range("A1").Value = "objectName"
$ws.Range("B1").Value = "actionType"
$ws.Range("C1").Value = "objectText"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.ColorIndex = 1

# Data rows
$ws.Range("A2").Value = "mns_File"
$ws.Range("C2").Value = "File"

$ws.Range("A3").Value = "mns_Settings"
$ws.Range("C3").Value = "Settings"

$ws.Range("A4").Value = "mns_Help"
$ws.Range("C4").Value = "Help"

# Column widths (best-fit like the sibling sheets)
$ws.Columns.Item(1).ColumnWidth = 10.05
$ws.Columns.Item(2).ColumnWidth = 9.05
$ws.Columns.Item(3).ColumnWidth = 8.72

# Freeze header row and set selection to A2
$selResult = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ------------------------------------------------------------------
# 2. messageBox sheet: add new row for the settings-language message
# ------------------------------------------------------------------
$mb = $wb.Worksheets.Item("messageBox")
$mb.Range("A43").Value = "mbx_frm_Settings_cbx_Language_TextChanged"
$mb.Range("C43").Value = "Please restart the app for the language change to take place."
$mb.Range("C43").WrapText = $true

$mb.Select()
$selResult = $mb.Range("A43").Select()

# ------------------------------------------------------------------
# 3. ToolStrip sheet: selection moved to C1
# ------------------------------------------------------------------
$ts = $wb.Worksheets.Item("ToolStrip")
$selResult = $ts.Range("C1").Select()

# ------------------------------------------------------------------
# Leave messageBox as the active sheet/tab (unchanged from before)
# ------------------------------------------------------------------
$mb.Select()

Write-Output "done"
